$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the Heading1
#    title paragraph ("Play Day and Night Free Slot Game - Review 2021").
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(2)
$p2.Style = "Normal"

$insertStart = $p2.Range.Start
$boldLabel = "Meta description"
$fullText = $boldLabel + ": Play Day and Night, an ancient Egyptian-themed slot game with stunning graphics and free spins. Read our review to learn more and play for free."

$insRange = $d.Range($insertStart, $insertStart)
$insRange.InsertAfter($fullText)

# Bold only the "Meta description" label, leaving the rest as plain text.
$boldEnd = $insertStart + $boldLabel.Length
$boldRange = $d.Range($insertStart, $boldEnd)
$boldRange.Bold = 1

# ---------------------------------------------------------------------
# 2) Remove the duplicated bold "Play Day and Night Free Slot Game -
#    Review 2021" paragraph that used to sit right before the closing
#    italic meta-description paragraph near the end of the document.
# ---------------------------------------------------------------------
$n = $d.Paragraphs.Count
for ($i = $n; $i -ge 1; $i--) {
    $cand = $d.Paragraphs($i)
    $candRange = $cand.Range
    if ($candRange.Text -eq "Play Day and Night Free Slot Game - Review 2021" -and $candRange.Bold -eq 1) {
        $candRange.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new
#    feature-image prompt text.
# ---------------------------------------------------------------------
$oldTail = "Play Day and Night, an ancient Egyptian-themed slot game with stunning graphics and free spins. Read our review to learn more and play for free."
$newTail = "Create a feature image fitting the game " + [char]34 + "Day and Night" + [char]34 + ": Please create an image in cartoon style featuring a happy Maya warrior with glasses, surrounded by the opposing forces of light and darkness. The warrior should be holding a staff or other ancient Egyptian weapon, and standing confidently amidst the clash of the two deities, Ra and Bastet. The image should be colorful and dynamic, with radiant orange hues on one side and a dark, mystical blue on the other. The name of the game, " + [char]34 + "Day and Night" + [char]34 + " should be featured prominently in the image, and any additional text or graphics should be inspired by ancient Egyptian mythology. Overall, the image should capture the excitement and intrigue of this unique and visually stunning online slot game."

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$findScope = $lastRange.Find
$findScope.ClearFormatting()
$findScope.Text = $oldTail
$findScope.Replacement.ClearFormatting()
$findScope.Replacement.Text = $newTail
$findScope.Execute($findScope.Text, $true, $false, $false, $false, $false, $true, 1, $false, $findScope.Replacement.Text, 2)
